# Add "Introduction" / "Content Panes" / "References" headings to the
# report: the existing paragraph becomes the bold "Introduction" heading,
# and three more bold paragraphs (one blank spacer, "Content Panes", one
# blank spacer, "References") are appended after it.

$d = $word.ActiveDocument

# --- Discover the existing formatting so we can re-use it verbatim ------
# Pull the raw package XML once so we can recover (a) the first paragraph's
# identity attributes (w14:paraId/w:rsid*) and (b) the run-level rPr that is
# already applied (rFonts/lang), so the new heading paragraphs stay in sync
# with whatever formatting this document already carries instead of us
# hard-coding it.
$full = $d.WordOpenXML

$p1Attrs = ""
if ($full -match '<w:p\s+([^>]*)>') {
    $p1Attrs = " " + $matches[1]
}

$baseRPr = '<w:lang w:val="en-US"/>'
if ($full -match '<w:body>.*?<w:r>\s*<w:rPr>(.*?)</w:rPr>') {
    $baseRPr = $matches[1]
}

# Build the bold version of that run-properties blob: insert <w:b/><w:bCs/>
# right after <w:rFonts .../> if present, otherwise just prepend it.
if ($baseRPr -match '(<w:rFonts[^/]*/>)') {
    $boldRPrInner = $baseRPr -replace [regex]::Escape($matches[1]), ($matches[1] + '<w:b/><w:bCs/>')
} else {
    $boldRPrInner = '<w:b/><w:bCs/>' + $baseRPr
}
$boldRPr = '<w:rPr>' + $boldRPrInner + '</w:rPr>'

# --- Helper to build a single bold heading paragraph ---------------------
function New-HeadingParagraphXml($text, $attrs) {
    if ($null -eq $text) {
        return "<w:p$attrs><w:pPr>$boldRPr</w:pPr></w:p>"
    }
    return "<w:p$attrs><w:pPr>$boldRPr</w:pPr><w:r>$boldRPr<w:t>$text</w:t></w:r></w:p>"
}

$newBody = New-HeadingParagraphXml "Introduction" $p1Attrs
$newBody += New-HeadingParagraphXml $null ""
$newBody += New-HeadingParagraphXml "Content Panes" ""
$newBody += New-HeadingParagraphXml $null ""
$newBody += New-HeadingParagraphXml "References" ""

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
    '<w:body>' + $newBody + '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

# Replace the (single) original paragraph's range with the five new
# paragraphs. The range spans the whole paragraph, including its paragraph
# mark, so InsertXML swaps the paragraph out for our replacement content
# while leaving the rest of the document (here, just the sectPr) untouched.
$targetRange = $d.Paragraphs(1).Range
$targetRange.InsertXML($packageXml)
